$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2566.4443
$ws.Range("I41").Value = 2871.1428
$ws.Range("K41").Value = 2871.1428
$ws.Range("M41").Value = -2431.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2659
$ws.Range("I86").Value = 966.6667
$ws.Range("K86").Value = 966.6667
$ws.Range("M86").Value = 156.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2659
$ws.Range("I89").Value = 966.6667
$ws.Range("K89").Value = 4833.3335
$ws.Range("M89").Value = 782.6665000000003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 21999
$ws.Range("J95").Value = 21999
$ws.Range("L95").Value = 21999
$ws.Range("N95").Value = -27491

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4499.5
$ws.Range("I113").Value = 3999
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 3999
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -745
$ws.Range("N113").Value = -11508

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 190
$ws.Range("I118").Value = 190
$ws.Range("K118").Value = 570
$ws.Range("M118").Value = 1087

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6590.9287
$ws.Range("I132").Value = 5427.3
$ws.Range("K132").Value = 16281.9
$ws.Range("M132").Value = -13751.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3976.5789
$ws.Range("I138").Value = 3893.3333
$ws.Range("J138").Value = 3992.1875
$ws.Range("K138").Value = 11679.9999
$ws.Range("L138").Value = 11976.5625
$ws.Range("M138").Value = -6539.999899999999
$ws.Range("N138").Value = -22256.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5322.6
$ws.Range("I74").Value = 3299.5
$ws.Range("K74").Value = 3299.5
$ws.Range("M74").Value = -2425.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5322.6
$ws.Range("I77").Value = 3299.5
$ws.Range("K77").Value = 16497.5
$ws.Range("M77").Value = -12129.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 4000
$ws.Range("K20").Value = 4000
$ws.Range("M20").Value = -3753

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2229
$ws.Range("I86").Value = 2229
$ws.Range("K86").Value = 2229
$ws.Range("M86").Value = -1106

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2229
$ws.Range("I89").Value = 2229
$ws.Range("K89").Value = 11145
$ws.Range("M89").Value = -5529

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1636.6666
$ws.Range("I94").Value = 1636.6666
$ws.Range("K94").Value = 1636.6666
$ws.Range("M94").Value = -1185.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4615.8823
$ws.Range("I134").Value = 770.9091
$ws.Range("J134").Value = 11665
$ws.Range("K134").Value = 2312.7273
$ws.Range("L134").Value = 34995
$ws.Range("M134").Value = 222.2727
$ws.Range("N134").Value = -40065

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 456
$ws.Range("I58").Value = 456
$ws.Range("K58").Value = 456
$ws.Range("M58").Value = -253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4199.1665
$ws.Range("I122").Value = 4049
$ws.Range("J122").Value = 4499.5
$ws.Range("K122").Value = 12147
$ws.Range("L122").Value = 13498.5
$ws.Range("M122").Value = -9697
$ws.Range("N122").Value = -18398.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6355.5
$ws.Range("I132").Value = 6355.5
$ws.Range("K132").Value = 19066.5
$ws.Range("M132").Value = -16536.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 456
$ws.Range("I136").Value = 456
$ws.Range("K136").Value = 1368
$ws.Range("M136").Value = 1182

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 212.75
$ws.Range("I2").Value = 216.33333
$ws.Range("J2").Value = 202
$ws.Range("K2").Value = 1297.99998
$ws.Range("L2").Value = 1212
$ws.Range("M2").Value = -1184.99998
$ws.Range("N2").Value = -1438

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1667023.4
$ws.Range("I4").Value = 1667023.4
$ws.Range("K4").Value = 5001070.199999999
$ws.Range("M4").Value = -5000958.199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 93.59999999999999
$ws.Range("I38").Value = 82.666664
$ws.Range("J38").Value = 110
$ws.Range("K38").Value = 247.999992
$ws.Range("L38").Value = 330
$ws.Range("M38").Value = 99.00000800000001
$ws.Range("N38").Value = -1024

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I68").Value = 1575
$ws.Range("J68").Value = 1075.6666
$ws.Range("K68").Value = 4725
$ws.Range("L68").Value = 3226.9998
$ws.Range("M68").Value = -3914
$ws.Range("N68").Value = -4848.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I71").Value = 1575
$ws.Range("J71").Value = 1075.6666
$ws.Range("K71").Value = 14175
$ws.Range("L71").Value = 9680.999400000001
$ws.Range("M71").Value = -10119
$ws.Range("N71").Value = -17792.9994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1300
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 1300
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2711
$ws.Range("J131").Value = 2711
$ws.Range("L131").Value = 8133
$ws.Range("N131").Value = -18213

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 5000
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 5000
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -4741
$ws.Range("N52").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2142.7144
$ws.Range("I97").Value = 1499.8334
$ws.Range("J97").Value = 6000
$ws.Range("K97").Value = 1499.8334
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -1003.8334
$ws.Range("N97").Value = -6992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 808.55554
$ws.Range("I102").Value = 808.55554
$ws.Range("K102").Value = 808.55554
$ws.Range("M102").Value = 813.44446

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 875.25
$ws.Range("I122").Value = 750.4
$ws.Range("K122").Value = 2251.2
$ws.Range("M122").Value = 198.8000000000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1675.8334
$ws.Range("I132").Value = 1309
$ws.Range("J132").Value = 2409.5
$ws.Range("K132").Value = 3927
$ws.Range("L132").Value = 7228.5
$ws.Range("M132").Value = -1397
$ws.Range("N132").Value = -12288.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6500
$ws.Range("I7").Value = 6500
$ws.Range("K7").Value = 6500
$ws.Range("M7").Value = -6388

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 749.6
$ws.Range("I16").Value = 749.6
$ws.Range("K16").Value = 749.6
$ws.Range("M16").Value = -579.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 3250.1
$ws.Range("I55").Value = 300.66666
$ws.Range("J55").Value = 4514.143
$ws.Range("K55").Value = 300.66666
$ws.Range("L55").Value = 4514.143
$ws.Range("M55").Value = -127.66666
$ws.Range("N55").Value = -4860.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1762.5
$ws.Range("I82").Value = 1783.3334
$ws.Range("K82").Value = 1783.3334
$ws.Range("M82").Value = -1422.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1762.5
$ws.Range("I85").Value = 1783.3334
$ws.Range("K85").Value = 1783.3334
$ws.Range("M85").Value = -535.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 8669.286
$ws.Range("I93").Value = 8669.286
$ws.Range("K93").Value = 8669.286
$ws.Range("M93").Value = -7421.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3410.8333
$ws.Range("I100").Value = 3393
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 3393
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -2852
$ws.Range("N100").Value = -4582

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4928.143
$ws.Range("I122").Value = 4999.25
$ws.Range("K122").Value = 14997.75
$ws.Range("M122").Value = -12547.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6500
$ws.Range("I126").Value = 6500
$ws.Range("K126").Value = 19500
$ws.Range("M126").Value = -17030

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5189.5
$ws.Range("I132").Value = 5342.857
$ws.Range("K132").Value = 16028.571
$ws.Range("M132").Value = -13498.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9575
$ws.Range("I136").Value = 11100
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 33300
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -30750
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 39559
$ws.Range("J45").Value = 41662
$ws.Range("L45").Value = 41662
$ws.Range("N45").Value = -42644

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 32500.092
$ws.Range("I81").Value = 35500.1
$ws.Range("K81").Value = 71000.2
$ws.Range("M81").Value = -69939.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 32500.092
$ws.Range("I84").Value = 35500.1
$ws.Range("K84").Value = 355001
$ws.Range("M84").Value = -349697

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3528.9473
$ws.Range("I132").Value = 2693.75
$ws.Range("J132").Value = 7983.3335
$ws.Range("K132").Value = 8081.25
$ws.Range("L132").Value = 23950.0005
$ws.Range("M132").Value = -5551.25
$ws.Range("N132").Value = -29010.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1200
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 3600
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -8700
